# S32K3xx fault map update
# - Rename device-family worksheets to reflect the broader set of covered parts
# - Add a new "Lockup and Lockstep errors" entry (with a yellow highlight) for
#   the C3 cell on the sheets that previously just said "Lockup"
# - Update the active sheet / selected cell bookkeeping to match the saved
#   workbook view state

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets -------------------------------------------------
$wsTitle   = $wb.Worksheets.Item(1)
$ws312     = $wb.Worksheets.Item(2)
$ws342     = $wb.Worksheets.Item(3)
$ws344     = $wb.Worksheets.Item(4)
$ws358     = $wb.Worksheets.Item(5)
$ws388     = $wb.Worksheets.Item(6)

$ws312.Name = "S32K310_S32K311_S32K312"
$ws342.Name = "S32K322_S32K341_S32K342"
$ws344.Name = "S32K314_S32K324_S32K344"
$ws358.Name = "S32K328_S32K338_S32K348_S32K358"

# --- 2. Update the "Lockup" cell on the affected sheets -------------------
# Title and S32K310_S32K311_S32K312 keep the plain "Lockup" text; the other
# four sheets get the new wording plus a yellow fill to flag the change.
# Build up the new formatting on the first sheet, then copy its formatting
# (format-painter style) onto the matching cells of the remaining sheets so
# that they all end up sharing the exact same cell style.
$c342 = $ws342.Range("C3")
$c342.Locked = $true
$c342.Interior.Color = 65535
$c342.Value = "Lockup and Lockstep errors"

$c342.Copy()

foreach ($ws in @($ws344, $ws358, $ws388)) {
    $c = $ws.Range("C3")
    $c.Value = "Lockup and Lockstep errors"
    $c.PasteSpecial(-4122)
}

# --- 3. Restore view state (active sheet / selected cells) ----------------
# Selecting a range activates its sheet, so set the sheets that should NOT
# stay active first, finishing with the one that should end up selected.
$wsTitle.Activate()
$wsTitle.Range("K6").Select()

$ws342.Activate()
$ws342.Range("C3").Select()

$ws344.Activate()
$ws344.Range("C3").Select()

$ws358.Activate()
$ws358.Range("C3").Select()

$ws388.Activate()
$ws388.Range("C3").Select()

$ws312.Activate()
$ws312.Range("C3").Select()
